$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title/caption in A1: batches 200 -> 64, confidence 99.00% -> 95.00%
$ws.Range("A1").Value = "Based on a simulation splitted into 64 batches and with 95.00% confidence"

# ---- NODE 1 block ----
$ws.Range("B4").Value = 0.52689600000000003
$ws.Range("D4").Value = 0.005398
$ws.Range("B5").Value = 2.7945479999999998
$ws.Range("D5").Value = 0.121612
$ws.Range("B6").Value = 0.79142400000000002
$ws.Range("D6").Value = 0.109916
$ws.Range("B7").Value = 2.002745
$ws.Range("D7").Value = 0.020379999999999999
$ws.Range("B8").Value = 5.3084389999999999
$ws.Range("D8").Value = 0.238816
$ws.Range("B9").Value = 1.5061990000000001
$ws.Range("D9").Value = 0.21133099999999999
$ws.Range("B10").Value = 0.76031000000000004
$ws.Range("D10").Value = 0.0085869999999999991

# ---- NODE 2 block ----
$ws.Range("B14").Value = 0.84731800000000002
$ws.Range("D14").Value = 0.011624000000000001
$ws.Range("B15").Value = 3.509001
$ws.Range("D15").Value = 0.072206000000000006
$ws.Range("B16").Value = 0.32366200000000001
$ws.Range("D16").Value = 0.047802999999999998
$ws.Range("B17").Value = 3.185117
$ws.Range("D17").Value = 0.037817999999999997
$ws.Range("B18").Value = 4.1482320000000001
$ws.Range("D18").Value = 0.099894999999999998
$ws.Range("B19").Value = 0.38419700000000001
$ws.Range("D19").Value = 0.057833000000000002
$ws.Range("B20").Value = 0.62729999999999997
$ws.Range("D20").Value = 0.0097439999999999992

# ---- NODE 3 block ----
$ws.Range("B24").Value = 1.3249979999999999
$ws.Range("D24").Value = 0.022620000000000001
$ws.Range("B25").Value = 3.3852669999999998
$ws.Range("D25").Value = 0.159745
$ws.Range("B26").Value = 0.89924499999999996
$ws.Range("D26").Value = 0.13547799999999999
$ws.Range("B27").Value = 2.4851749999999999
$ws.Range("D27").Value = 0.039510000000000003
$ws.Range("B28").Value = 2.5686010000000001
$ws.Range("D28").Value = 0.13736499999999999
$ws.Range("B29").Value = 0.68695899999999999
$ws.Range("D29").Value = 0.10886800000000001
$ws.Range("B30").Value = 0.62709599999999999
$ws.Range("D30").Value = 0.013576

# ---- NODE 4 block ----
$ws.Range("B34").Value = 0.386467
$ws.Range("D34").Value = 0.0029120000000000001
$ws.Range("B35").Value = 2.119929
$ws.Range("D35").Value = 0.045718000000000002
$ws.Range("B36").Value = 0.82004600000000005
$ws.Range("D36").Value = 0.036830000000000002
$ws.Range("B37").Value = 1.299774
$ws.Range("D37").Value = 0.012526000000000001
$ws.Range("B38").Value = 5.479495
$ws.Range("D38").Value = 0.107525
$ws.Range("B39").Value = 2.117353
$ws.Range("D39").Value = 0.089992000000000003
$ws.Range("B40").Value = 0.84048400000000001
$ws.Range("D40").Value = 0.0084239999999999992
$ws.Range("B41").Value = 0.041488999999999998
$ws.Range("D41").Value = 0.0040829999999999998

# ---- Summary row (43) ----
$ws.Range("B43").Value = 11.8087
$ws.Range("D43").Value = 0.24

# Update the selection to match the new view (A1:F1, the title row)
$ws.Range("A1:F1").Select()
